$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.415.43"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.572.91"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.90"
$ws.Range("E5").Value = "  +3.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.67"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "3.567.94"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.618"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.216"
$ws.Range("E10").Value = "  +7.73%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.02"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.51"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "4.136.45"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "70.471.22"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.76"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.568.18"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.05"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "580.36"
$ws.Range("E20").Value = "  +5.14%  "
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.39"
$ws.Range("E23").Value = "  -2.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.72"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.65"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.96"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.42"
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.34"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.08"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.28"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.53"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.71"
$ws.Range("E35").Value = "  +20.22%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "533.37"
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.405"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.37"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("D42").Value = "3.523.48"
$ws.Range("E42").Value = "  +4.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.53"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0461"
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.141"
$ws.Range("E48").Value = "  +3.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.25"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.55"
$ws.Range("E51").Value = "  -0.64%  "
